$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new date column header (BY1) - mirrors the existing "25-sep" column (BX)
$ws.Range("BY1").Value = "26-sep"

# Fill in the new BY column values for each data row
$ws.Range("BY2").Value  = 10
$ws.Range("BY3").Value  = 14
$ws.Range("BY4").Value  = 11
$ws.Range("BY5").Value  = 11
$ws.Range("BY6").Value  = 8
$ws.Range("BY7").Value  = 13
$ws.Range("BY8").Value  = 15
$ws.Range("BY9").Value  = 15
$ws.Range("BY10").Value = 13
$ws.Range("BY11").Value = 8

# Match the centered, 0-decimal number style used by the rest of the table
$dataRange = $ws.Range("BY2:BY11")
$dataRange.HorizontalAlignment = -4108
$dataRange.NumberFormat = "0"

# Update the view: scroll/selection moves to BY8, no frozen top-left offset
$ws.Range("BY8").Select() | Out-Null
